# "New forms for the BMGF demo"
# Adds two new ODK "external_link" form entries (visit, plot) to the
# survey sheet's form list and to the choices sheet's test_forms list.

$wb = $excel.ActiveWorkbook
$wsSurvey  = $wb.Worksheets.Item("survey")
$wsChoices = $wb.Worksheets.Item("choices")

# ---------------------------------------------------------------------
# survey sheet: append the two new 3-row form blocks (rows 53-58),
# mirroring the existing "external_link" block pattern used by every
# other form (see rows 47-52 for "adult_coverage").
# ---------------------------------------------------------------------

function Add-FormBlock($ws, $startRow, $formName, $hashFormula) {
    $rTitle = $startRow
    $rLink  = $startRow + 1
    $rExit  = $startRow + 2

    # Title row: "<formName>"
    $ws.Rows.Item($rTitle).RowHeight = 17.5
    $ws.Cells.Item($rTitle, 1).Value = $formName

    # Link row: hash formula string (column B, quote-prefixed text),
    # "url" (E) and "Open form" (G) labels.
    $ws.Rows.Item($rLink).RowHeight = 66
    $ws.Cells.Item($rLink, 1).Value = ""
    $ws.Cells.Item($rLink, 1).Style = $ws.Cells.Item(48, 1).Style
    $ws.Cells.Item($rLink, 2).Value = $hashFormula
    $ws.Cells.Item($rLink, 5).Value = "external_link"
    $ws.Cells.Item($rLink, 7).Value = "Open form"

    # Exit-section row.
    $ws.Rows.Item($rExit).RowHeight = 17
    $ws.Cells.Item($rExit, 1).Value = ""
    $ws.Cells.Item($rExit, 1).Style = $ws.Cells.Item(48, 1).Style
    $ws.Cells.Item($rExit, 2).Value = ""
    $ws.Cells.Item($rExit, 2).Style = $ws.Cells.Item(48, 1).Style
    $ws.Cells.Item($rExit, 3).Value = "exit section"
}

# NB: the leading "'" here is intentionally doubled. A single leading
# apostrophe is consumed by Excel's quote-prefix ("treat as text") input
# convention and stripped from the stored value; doubling it yields a
# stored string whose first character genuinely is "'" (matching the
# existing sibling rows, e.g. row 48) while still tripping the
# quote-prefix formatting (style s="4") that those rows use.
Add-FormBlock $wsSurvey 53 "visit" "''?' + opendatakit.getHashString('../tables/visit/forms/visit/',null)"
Add-FormBlock $wsSurvey 56 "plot"  "''?' + opendatakit.getHashString('../tables/plot/forms/plot/',null)"

$wsSurvey.Range("B58").Select()

# ---------------------------------------------------------------------
# choices sheet: two new "test_forms" choice rows (rows 17-18).
# ---------------------------------------------------------------------

function Add-ChoiceRow($ws, $row, $dataValue, $displayText) {
    $ws.Cells.Item($row, 1).Value = "test_forms"
    $ws.Cells.Item($row, 1).Style = $ws.Cells.Item(16, 1).Style
    $ws.Cells.Item($row, 2).Value = $dataValue
    $ws.Cells.Item($row, 3).Value = $displayText
}

Add-ChoiceRow $wsChoices 17 "visit" "Visit"
Add-ChoiceRow $wsChoices 18 "plot"  "Plot"

$wsChoices.Range("B19").Select()
$wsChoices.Activate()
